$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-25 Sunday" "2024-02-26 Monday"

Replace-Text "19×56=1064" "80×46=3680"
Replace-Text "89×99=8811" "22×72=1584"
Replace-Text "74×76=5624" "67×87=5829"
Replace-Text "63×29=1827" "44×38=1672"
Replace-Text "15×29=435" "53×73=3869"

Replace-Text "41×53=2173" "26×40=1040"
Replace-Text "23×88=2024" "40×89=3560"
Replace-Text "24×50=1200" "94×38=3572"
Replace-Text "12×81=972" "80×73=5840"
Replace-Text "98×25=2450" "11×65=715"

Replace-Text "99×27=2673" "80×38=3040"
Replace-Text "20×87=1740" "82×59=4838"
Replace-Text "86×56=4816" "41×59=2419"
Replace-Text "57×46=2622" "65×27=1755"
Replace-Text "34×23=782" "73×61=4453"

Replace-Text "68×28=1904" "98×60=5880"
Replace-Text "95×82=7790" "71×18=1278"
Replace-Text "24×21=504" "84×55=4620"
Replace-Text "43×53=2279" "32×89=2848"
Replace-Text "63×25=1575" "36×23=828"

Replace-Text "89×33=2937" "26×34=884"
Replace-Text "36×86=3096" "47×43=2021"
Replace-Text "52×69=3588" "93×67=6231"
Replace-Text "66×71=4686" "58×37=2146"
Replace-Text "55×49=2695" "99×14=1386"
